$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Columns A,B,C,D,E,F,G,H,J,K,L,O,P are stored as text in the source
# workbook, even when the content looks numeric/date-like (e.g. the case
# number, comuna number, OT number or the date string). Force the number
# format to Text on those cells before assigning so Excel does not silently
# convert them to numbers/dates, then restore the default "Normal" style so
# no extra formatting is left behind on the cell (matching the rest of the
# sheet's plain data cells).
$textCols = 1,2,3,4,5,6,7,8,10,11,12,15,16
foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "4088"
$ws.Cells.Item($row, 2).Value = "9/2/2025"
$ws.Cells.Item($row, 3).Value = "GOMEZ, VALENTIN 3648"
$ws.Cells.Item($row, 4).Value = "5"
$ws.Cells.Item($row, 5).Value = "809406164"
$ws.Cells.Item($row, 6).Value = "AYKO"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Picada"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.416416
$ws.Cells.Item($row, 14).Value = -34.604812
$ws.Cells.Item($row, 15).Value = "Almagro"
$ws.Cells.Item($row, 16).Value = "Capital Sur"

foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).Style = "Normal"
}
